# Connect to database and initial analysis.
# Performance sheet gains a "Type" column (breaking the Documents/Time
# counts out per dataset) and three new rows for the extra datasets that
# were queried (product, product reviews, price overall).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance")

# Make room for the new leading "Type" column; the existing Documents /
# Time in seconds / Time in minutes columns (and the minutes formula that
# references them) shift right automatically.
$ws.Columns.Item(1).Insert() | Out-Null

# Header
$ws.Range("A1").Value = "Type"
$ws.Range("A1").Font.Bold = $true

# Row 2 stays "metadata" (its Documents/seconds/minutes numbers already
# existed before the column insert).
$ws.Range("A2").Value = "metadata"

# Query the remaining datasets, in the order they were looked up.
$ws.Range("A5").Value = "product reviews"
$ws.Range("B5").Value = 2952306

$ws.Range("A3").Value = "review"
$ws.Range("B3").Value = 17615000

$ws.Range("A4").Value = "product"
$ws.Range("B4").Value = 16964379

$ws.Range("A6").Value = "price overall"
$ws.Range("B6").Value = 2467052

# Resize the new column to fit its contents.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave focus on the Performance tab, on the next empty row.
$ws.Activate() | Out-Null
$ws.Range("B8").Select() | Out-Null
